# Clears the fixed "boilerplate" title/body text that was being
# displayed instead of the user's table input, leaving a single blank
# space in each affected run (fixes bug where static sample text showed
# up instead of the user-entered content).

$d = $word.ActiveDocument

$wdReplaceAll = 2
$wdFindContinue = 1

# 1) "Test Evaluation Report – Parte fixa" -> " "
$d.Content.Find.Execute(
    "Test Evaluation Report – Parte fixa",
    $true, $false, $false, $false, $false, $true,
    $wdFindContinue, $false, " ", $wdReplaceAll) | Out-Null

# 2) & 3) the two sample "Lorem ipsum ..." paragraphs -> " " (ReplaceAll
#    also happens to pick up the third Lorem ipsum paragraph further down,
#    which the diff shows is cleared too)
$d.Content.Find.Execute(
    "Lorem ipsum dolor sit amet, consectetur adipiscing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum.",
    $true, $false, $false, $false, $false, $true,
    $wdFindContinue, $false, " ", $wdReplaceAll) | Out-Null

# 4) "Test case execution summary" -> " "
$d.Content.Find.Execute(
    "Test case execution summary",
    $true, $false, $false, $false, $false, $true,
    $wdFindContinue, $false, " ", $wdReplaceAll) | Out-Null

# 5) "Parte fixa open risks and mitigation" -> " "
$d.Content.Find.Execute(
    "Parte fixa open risks and mitigation",
    $true, $false, $false, $false, $false, $true,
    $wdFindContinue, $false, " ", $wdReplaceAll) | Out-Null
